$wb = $excel.ActiveWorkbook
$nl = [char]10

# 1. Wrap the long land-use labels onto multiple lines (Alt+Enter line breaks)
#    so the table headers/labels read on two or three lines instead of one.
$replacements = @{
    "Impervious Roads"            = "Impervious" + $nl + "Roads";
    "Impervious Non-Roads"        = "Impervious" + $nl + "Non-Roads";
    "Tree Canopy Over Impervious" = "Tree Canopy" + $nl + "Over" + $nl + "Impervious";
    "Tree Canopy over Turf Grass" = "Tree Canopy" + $nl + "over Turf" + $nl + "Grass";
    "Wetlands (Other)"            = "Wetlands" + $nl + "(Other)";
    "Wetlands (Floodplain)"       = "Wetlands" + $nl + "(Floodplain)";
}

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 2)
        $val = $cell.Value()
        if ($replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}

# 2. On the "region" sheet, swap the "Tree Canopy Over Impervious" /
#    "Tree Canopy over Turf Grass" rows within each region block so the
#    Turf-Grass-canopy row now comes first.
$region = $wb.Worksheets.Item("region")
$rowPairs = @(8, 19, 30, 41, 52, 63)
foreach ($r1 in $rowPairs) {
    $r2 = $r1 + 1

    $b1 = $region.Cells.Item($r1, 2).Value()
    $c1 = $region.Cells.Item($r1, 3).Value()
    $b2 = $region.Cells.Item($r2, 2).Value()
    $c2 = $region.Cells.Item($r2, 3).Value()

    $region.Cells.Item($r1, 2).Value = $b2
    $region.Cells.Item($r1, 3).Value = $c2
    $region.Cells.Item($r2, 2).Value = $b1
    $region.Cells.Item($r2, 3).Value = $c1
}
